$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (file_name): Ctrl_plate1_leukocytes_test.fcs -> Ctrl_plate0_test.fcs,
#                        P15_D1_leukocytes_test.fcs -> P99_L1_test.fcs
$ws.Range("A2").Value = "Ctrl_plate0_test.fcs"
$ws.Range("A3").Value = "P99_L1_test.fcs"

# Column C (sample_id): Ctrl_plate1 -> Ctrl_plate0, P15_D1 -> P99_L1
$ws.Range("C2").Value = "Ctrl_plate0"
$ws.Range("C3").Value = "P99_L1"

# Column D (anchor): Ctrl stays Ctrl, P15_D1 -> P99_L1
$ws.Range("D2").Value = "Ctrl"
$ws.Range("D3").Value = "P99_L1"

# Column B (run): 1 -> 2, 5 -> 6
$ws.Range("B2").Value = 2
$ws.Range("B3").Value = 6

# Update the active selection to D6 (was D9)
$ws.Range("D6").Select()
